# Aggiornamento dati Maranello fino al 9 agosto 2021
# Appends daily rows (329-343) for date serials 44403..44417, mirroring the
# layout (columns A=data, B=nuovi pos., C=somma mobile 7gg., D=somma mobile
# 7gg. per 100mila abitanti) and date-cell formatting of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing date cell (A328, style "s=2")
# onto the new date cells A329:A343 before writing values, so the new rows
# keep the same centered/bordered date-number format as the rest of column A.
$ws.Range("A328").Copy()
$ws.Range("A329:A343").PasteSpecial(-4122)

$newRows = @(
    @(329, 44403, 1, 10, 56.91843588138198),
    @(330, 44404, 0, 10, 56.91843588138198),
    @(331, 44405, 0, 7, 39.84290511696739),
    @(332, 44406, 1, 5, 28.45921794069099),
    @(333, 44407, 2, 7, 39.84290511696739),
    @(334, 44408, 2, 7, 39.84290511696739),
    @(335, 44409, 2, 8, 45.53474870510559),
    @(336, 44410, 1, 8, 45.53474870510559),
    @(337, 44411, 0, 8, 45.53474870510559),
    @(338, 44412, 2, 10, 56.91843588138198),
    @(339, 44413, 0, 9, 51.22659229324378),
    @(340, 44414, 3, 10, 56.91843588138198),
    @(341, 44415, 4, 12, 68.30212305765838),
    @(342, 44416, 2, 12, 68.30212305765838),
    @(343, 44417, 0, 11, 62.61027946952018)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
